$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 38, shifting existing rows 38:60 down to 39:61
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record
$ws.Range("A38").Value = 9
$ws.Range("B38").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C38").Value = "Metropolitana"
$ws.Range("D38").Value = 45062
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = 100112010
$ws.Range("G38").Value = "Achicoria"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 52
$ws.Range("K38").Value = 7000
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = 7000
$ws.Range("N38").Value = "$/caja 16 unidades"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 438
$ws.Range("Q38").Value = 16
$ws.Range("R38").Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D
$ws.Range("D38").NumberFormat = $ws.Range("D39").NumberFormat
